# Append new form-submission rows (40-53) to the sheet, matching the
# source diff. Columns: A First Name, B Last Name, C Email, D Company
# Email, E Phone Number, F Organization, G Help Request, H Budget,
# I Services.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("User","User","akshanaggarwal2005@hotmail.com","null","Not provided","Not provided","Interested in services","Not specified","General Inquiry"),
    @("Krissh","Gera","gkrissh7@gmail.com","3dcarcare@gmail.com","09870305778","3d car care India","make meta ads for me","2000","Website building"),
    @("Akshan","Aggarwal","akshanaggarwal2005@hotmail.com","3dcarcare@gmail.com","09870305778","3d car care India","make meta ads for me","2000","Website building"),
    @("Akshan","Aggarwal","akshanaggarwal2005@hotmail.com","3dcarcare@gmail.com","09870305778","3d car care India","make meta ads for me","2000","Website building"),
    @("Krishna","Aggarwal","aggarwalkrishna3163@gmail.com","3dcarcare@gmail.com","09870305778","3d car care India","make meta ads for me","2000","Website building"),
    @("User","User","official.pranav02@gmail.com","null","Not provided","Not provided","Interested in services","Not specified","General Inquiry"),
    @("Akshan","Aggarwal","akshanaggarwal20@gmail.com","3dcarcare@gmail.com","09870305778","3d car care India","make meta ads for me","2000","Website building"),
    @("User","User","akshanaggarwal2005@hotmail.com","null","Not provided","Not provided","Interested in services","Not specified","General Inquiry"),
    @("User","User","akshanaggarwal2005@hotmail.com","null","Not provided","Not provided","Interested in services","Not specified","General Inquiry"),
    @("User","User","akshanaggarwal2005@hotmail.com","null","Not provided","Not provided","Interested in services","Not specified","General Inquiry"),
    @("Akshan","Aggarwal","akshanaggarwal2005@hotmail.com","3dcarcare@gmail.com","09870305778","3d car care India","make meta ads for me","2000","Website building"),
    @("Akshan","Aggarwal","akshanaggarwal2005@hotmail.com","3dcarcare@gmail.com","09870305778","3d car care India","make meta ads for me","2000","Website building"),
    @("Akshan","Aggarwal","akshanaggarwal2005@hotmail.com","3D","9870305778","3D","No help","50000","Website development "),
    @("User","User","akshanaggarwal2005@hotmail.com","null","Not provided","Not provided","Interested in services","Not specified","General Inquiry")
)

$startRow = 40
$columns = @("A","B","C","D","E","F","G","H","I")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $cellRef = $columns[$j] + $r
        $cell = $ws.Range($cellRef)
        $value = $rowData[$j]
        # Force text storage only for numeric-looking strings (phone
        # numbers, budgets, etc.) so they keep their original form
        # (e.g. leading zeros) instead of being coerced to numbers.
        if ($value -match '^-?\d+(\.\d+)?$') {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $value
    }
}
